# issue #5: stock data output to json file
# Adds a new "property_category" column (value "stock") to the 股票 (stock)
# worksheet, positioned right before the existing "date" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (column H)
$ws.Columns("H:H").Insert()

# Header for the newly inserted column
$ws.Cells.Item(1, 8).Value = "property_category"

# Fill the value for every data row (rows 2-6) with "stock"
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
  $ws.Cells.Item($r, 8).Value = "stock"
}
